$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.952336311340332
$ws.Range("B1").Value = 7.339256763458252
$ws.Range("C1").Value = 4.07763147354126
$ws.Range("D1").Value = 2.085773468017578
$ws.Range("E1").Value = 1.44983446598053
